$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.180.14"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +0.96%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.901.71"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +1.43%  "

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.18%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "306.66"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +0.04%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.000"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +0.12%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5230"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +1.42%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3767"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.77%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07238"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +0.66%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "21.16"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +2.12%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.8978"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -0.18%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08427"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +11.60%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.903.17"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +1.55%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "94.54"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -0.40%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.263"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +0.13%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.001"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +0.18%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000008588"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +0.83%  "

$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +1.55%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.000"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +0.13%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "27.225.06"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +1.01%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.053"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +0.37%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.142.28"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +2.06%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.58"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +1.69%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.414"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -0.11%  "

$ws.Range("B25").NumberFormat = "@"
$ws.Range("B25").Value = "LidoDAOToken"
$ws.Range("C25").NumberFormat = "@"
$ws.Range("C25").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.285"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +8.20%  "

$ws.Range("B26").NumberFormat = "@"
$ws.Range("B26").Value = "Monero"
$ws.Range("C26").NumberFormat = "@"
$ws.Range("C26").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "146.65"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +0.36%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.753"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -1.50%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.14"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +0.83%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "114.77"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -0.16%  "

$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +0.19%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.788"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +0.28%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09211"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +0.28%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.8119"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +7.87%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.05054"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +0.52%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.238"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +5.64%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.961"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -0.97%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.366"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +3.54%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.564"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +3.16%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.5689"

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01977"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -0.53%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.071"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +0.01%  "

$ws.Range("B42").NumberFormat = "@"
$ws.Range("B42").Value = "FraxShare"
$ws.Range("C42").NumberFormat = "@"
$ws.Range("C42").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.637"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +0.86%  "

$ws.Range("B43").NumberFormat = "@"
$ws.Range("B43").Value = "Aptos"
$ws.Range("C43").NumberFormat = "@"
$ws.Range("C43").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.963"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +3.17%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "118.23"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +2.21%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.1510"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +0.67%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4821"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +0.96%  "

$ws.Range("B47").NumberFormat = "@"
$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").NumberFormat = "@"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "10.17"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +0.40%  "

$ws.Range("B48").NumberFormat = "@"
$ws.Range("B48").Value = "PaxDollar"
$ws.Range("C48").NumberFormat = "@"
$ws.Range("C48").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.000"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +0.13%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.610"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +2.87%  "

$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +0.91%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "63.57"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +0.28%  "
